# Suite à l'ajout des numéros de téléphone étrangers, les numéros de
# test précédemment utilisés dans le fichier d'exemple d'import usagers
# (format "01...") sont remplacés par un format "06..." pour la colonne
# "Numéro de téléphone" (H2:H20) de la feuille active ("Feuil1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$phones = @(
    "0601010101",
    "0601010102",
    "0601010103",
    "0601010104",
    "0601010105",
    "0601010106",
    "0601010107",
    "0601010108",
    "0601010109",
    "0601010110",
    "0601010111",
    "0601010112",
    "0601010113",
    "0601010114",
    "0601010115",
    "0601010116",
    "0601010117",
    "0601010118",
    "0601010119"
)

for ($i = 0; $i -lt $phones.Length; $i++) {
    $row = $i + 2
    $ws.Range("H" + $row).Value = $phones[$i]
}

# Reflect the selection left behind after editing the phone-number column.
$ws.Activate()
$ws.Range("H2:H20").Select()
